$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so numeric-looking strings
# (e.g. "214.84") are stored as text, matching the inlineStr cells
# already present in this workbook, not auto-converted to numbers.
$cells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "E6", "E7", "E8", "D9", "E9", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "E14", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "E19", "D20", "E20", "E21", "E22", "D23", "E23", "E24", "D25", "E27", "D28", "E28", "E29", "E30", "E31", "E32", "D33", "E33", "E34", "E35", "E36", "E37", "E38", "E39", "D40", "E40", "E41", "D42", "E42", "D43", "E43", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D49", "E49", "E50", "D51", "E51")
foreach ($addr in $cells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "26.985.90"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.671.25"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "214.84"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("D9").Value = "21.38"
$ws.Range("E9").Value = "  +5.23%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "0.0888"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "1.908.12"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "1.640.30"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "66.11"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "26.984.54"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("D19").Value = "234.75"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "0.0₃0735"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").Value = "9.25"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").Value = "147.01"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("D28").Value = "0.112"
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "1.535.20"
$ws.Range("E33").Value = "  +6.07%  "
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("D40").Value = "1.05"
$ws.Range("E40").Value = "  +4.91%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "67.55"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("D43").Value = "5.53"
$ws.Range("E43").Value = "  -3.52%  "
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("D45").Value = "1.815.06"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("D46").Value = "0.780"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "90.34"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D49").Value = "1.53"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("D51").Value = "8.01"
$ws.Range("E51").Value = "  +6.12%  "
